# Rebuild the "Test" worksheet:
#  - drop the "allowed" (C) and "dif" (D) columns entirely
#  - rename the "Bitcoin" asset to "CLCD16"
#  - reorder each fund's asset rows to: Stocks, LFT, CLCD16, PETR4

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "allowed" and "dif" columns, shifting the remaining cells left.
$ws.Range("C1:D13").Delete(-4159) | Out-Null

# Column B values, row by row (row 1 is the header and is unchanged).
$assetNames = @(
    "Stocks", "LFT", "CLCD16", "PETR4",
    "Stocks", "LFT", "CLCD16", "PETR4",
    "Stocks", "LFT", "CLCD16", "PETR4"
)

for ($i = 0; $i -lt $assetNames.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $assetNames[$i]
}

$ws.Range("B11").Select() | Out-Null
